# Smith Creek Deployment Log - add location names for the two waterLevel
# sensors (N 23rd St / N Kerr Ave) next to their coordinate rows, and
# normalize the header row's cell formatting so it matches column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (C2:H2) was carrying a redundant "applyAlignment" style that has no
# actual alignment set - functionally identical to B2's style. Re-apply B2's
# font across the row so these cells settle back onto the same style as B2.
$ws.Range("C2:H2").Font.Name = $ws.Range("B2").Font.Name()

# waterLevel_0000 (row 21) is at the N 23rd St bridge site; waterLevel_0001
# (row 22) is at the N Kerr Ave site. Record the site names in column C.
$ws.Range("C21").Value = "N 23rd St"
$ws.Range("C22").Value = "N Kerr Ave"

# Leave the saved selection on C23, as in the updated workbook view.
[void]$ws.Range("C23").Select()
